$d = $word.ActiveDocument

# Locate the "Impact" sub-heading inside "KEY ACHIEVEMENTS AND IMPACT" -
# the six achievement bullets immediately follow it. Anchoring on the
# heading (rather than a hard-coded paragraph number) keeps this robust
# even if unrelated paragraphs earlier in the doc shift.
$impactHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Impact" -and $p.Style.NameLocal -eq "Heading 3") {
        $impactHeadingIndex = $i
        break
    }
}

if ($impactHeadingIndex -eq -1) {
    throw "Could not find the 'Impact' sub-heading under KEY ACHIEVEMENTS AND IMPACT"
}

$b1 = $impactHeadingIndex + 1  # • Achieved 87% prediction accuracy...
$b2 = $impactHeadingIndex + 2  # • Delivered $4.9M additional revenue...
$b3 = $impactHeadingIndex + 3  # • Built redistricting platform...
$b4 = $impactHeadingIndex + 4  # • Developed longitudinal data analysis methods...
$b5 = $impactHeadingIndex + 5  # • Discovered systematic race coding errors...
$b6 = $impactHeadingIndex + 6  # • Trigonometric algorithm for boundary estimation...

# Delete the two trailing bullets first (highest index first) so the
# earlier paragraph indices we still need stay valid.
$d.Paragraphs.Item($b6).Range.Delete()
$d.Paragraphs.Item($b5).Range.Delete()

# Rewrite the remaining four bullets with the new impact-focused statements.
$d.Paragraphs.Item($b4).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$d.Paragraphs.Item($b3).Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
$d.Paragraphs.Item($b2).Range.Text = "• 23% conversion rate improvement"
$d.Paragraphs.Item($b1).Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
